$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3940
$ws.Cells.Item(40, 9).Value = 4675
$ws.Cells.Item(40, 10).Value = 1000
$ws.Cells.Item(40, 11).Value = 4675
$ws.Cells.Item(40, 12).Value = 1000
$ws.Cells.Item(40, 13).Value = -4500
$ws.Cells.Item(40, 14).Value = -1350

$ws.Cells.Item(94, 8).Value = 1483.7142
$ws.Cells.Item(94, 9).Value = 862
$ws.Cells.Item(94, 11).Value = 862
$ws.Cells.Item(94, 13).Value = -411

$ws.Cells.Item(111, 8).Value = 10016
$ws.Cells.Item(111, 9).Value = 10000
$ws.Cells.Item(111, 11).Value = 30000
$ws.Cells.Item(111, 13).Value = -26933

$ws.Cells.Item(141, 8).Value = 1873.75
$ws.Cells.Item(141, 9).Value = 1998.3334
$ws.Cells.Item(141, 11).Value = 5995.0002
$ws.Cells.Item(141, 13).Value = -815.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15215.154
$ws.Cells.Item(32, 9).Value = 15215.154
$ws.Cells.Item(32, 11).Value = 15215.154
$ws.Cells.Item(32, 13).Value = -14928.154

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(69, 8).Value = 40295
$ws.Cells.Item(69, 10).Value = 40295
$ws.Cells.Item(69, 12).Value = 40295
$ws.Cells.Item(69, 14).Value = -41917

$ws.Cells.Item(72, 8).Value = 40295
$ws.Cells.Item(72, 10).Value = 40295
$ws.Cells.Item(72, 12).Value = 120885
$ws.Cells.Item(72, 14).Value = -128997

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(39, 8).Value = 4404
$ws.Cells.Item(39, 10).Value = 5757
$ws.Cells.Item(39, 12).Value = 5757
$ws.Cells.Item(39, 14).Value = -6539

$ws.Cells.Item(49, 8).Value = 4404
$ws.Cells.Item(49, 10).Value = 5757
$ws.Cells.Item(49, 12).Value = 5757
$ws.Cells.Item(49, 14).Value = -6121

$ws.Cells.Item(58, 8).Value = 5449.75
$ws.Cells.Item(58, 9).Value = 2299.5
$ws.Cells.Item(58, 10).Value = 8600
$ws.Cells.Item(58, 11).Value = 2299.5
$ws.Cells.Item(58, 12).Value = 8600
$ws.Cells.Item(58, 13).Value = -2096.5
$ws.Cells.Item(58, 14).Value = -9006

$ws.Cells.Item(99, 8).Value = 6982.4
$ws.Cells.Item(99, 9).Value = 7478
$ws.Cells.Item(99, 11).Value = 7478
$ws.Cells.Item(99, 13).Value = -5980

$ws.Cells.Item(122, 8).Value = 3603.3333
$ws.Cells.Item(122, 9).Value = 2911
$ws.Cells.Item(122, 11).Value = 8733
$ws.Cells.Item(122, 13).Value = -6283

$ws.Cells.Item(125, 8).Value = 15000
$ws.Cells.Item(125, 10).Value = 15000
$ws.Cells.Item(125, 12).Value = 15000
$ws.Cells.Item(125, 14).Value = -19920

$ws.Cells.Item(126, 8).Value = 6982.4
$ws.Cells.Item(126, 9).Value = 7478
$ws.Cells.Item(126, 11).Value = 22434
$ws.Cells.Item(126, 13).Value = -19964

$ws.Cells.Item(132, 8).Value = 3932.6667
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 3932.6667
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 11798.0001
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -16858.0001

$ws.Cells.Item(136, 8).Value = 5449.75
$ws.Cells.Item(136, 9).Value = 2299.5
$ws.Cells.Item(136, 10).Value = 8600
$ws.Cells.Item(136, 11).Value = 6898.5
$ws.Cells.Item(136, 12).Value = 25800
$ws.Cells.Item(136, 13).Value = -4348.5
$ws.Cells.Item(136, 14).Value = -30900

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 856
$ws.Cells.Item(5, 9).Value = 1148.25
$ws.Cells.Item(5, 10).Value = 466.33334
$ws.Cells.Item(5, 11).Value = 3444.75
$ws.Cells.Item(5, 12).Value = 1399.00002
$ws.Cells.Item(5, 13).Value = -3332.75
$ws.Cells.Item(5, 14).Value = -1623.00002

$ws.Cells.Item(46, 8).Value = 1979
$ws.Cells.Item(46, 10).Value = 1979
$ws.Cells.Item(46, 12).Value = 5937
$ws.Cells.Item(46, 14).Value = -6119

$ws.Cells.Item(113, 8).Value = 1345
$ws.Cells.Item(113, 9).Value = 697.4
$ws.Cells.Item(113, 11).Value = 2092.2
$ws.Cells.Item(113, 13).Value = 77.80000000000018

$ws.Cells.Item(122, 8).Value = 3001.0667
$ws.Cells.Item(122, 9).Value = 1426
$ws.Cells.Item(122, 11).Value = 12834
$ws.Cells.Item(122, 13).Value = -10384

$ws.Cells.Item(128, 8).Value = 199999
$ws.Cells.Item(128, 9).Value = 199999
$ws.Cells.Item(128, 11).Value = 599997
$ws.Cells.Item(128, 13).Value = -595017

$ws.Cells.Item(135, 8).Value = 856
$ws.Cells.Item(135, 9).Value = 1148.25
$ws.Cells.Item(135, 10).Value = 466.33334
$ws.Cells.Item(135, 11).Value = 10334.25
$ws.Cells.Item(135, 12).Value = 4197.00006
$ws.Cells.Item(135, 13).Value = -7799.25
$ws.Cells.Item(135, 14).Value = -9267.00006

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2822.3
$ws.Cells.Item(7, 9).Value = 2802.7778
$ws.Cells.Item(7, 10).Value = 2998
$ws.Cells.Item(7, 11).Value = 2802.7778
$ws.Cells.Item(7, 12).Value = 2998
$ws.Cells.Item(7, 13).Value = -2690.7778
$ws.Cells.Item(7, 14).Value = -3222

$ws.Cells.Item(40, 8).Value = 4399.857
$ws.Cells.Item(40, 9).Value = 4249.75
$ws.Cells.Item(40, 11).Value = 4249.75
$ws.Cells.Item(40, 13).Value = -4113.75

$ws.Cells.Item(126, 8).Value = 2822.3
$ws.Cells.Item(126, 9).Value = 2802.7778
$ws.Cells.Item(126, 10).Value = 2998
$ws.Cells.Item(126, 11).Value = 8408.3334
$ws.Cells.Item(126, 12).Value = 8994
$ws.Cells.Item(126, 13).Value = -5938.3334
$ws.Cells.Item(126, 14).Value = -13934

$ws.Cells.Item(136, 8).Value = 3498.5
$ws.Cells.Item(136, 9).Value = 3331.3333
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 11).Value = 9993.999899999999
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = -7443.999899999999
$ws.Cells.Item(136, 14).Value = -17100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4055.6
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 4055.6
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 4055.6
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).Value = -5303.6

$ws.Cells.Item(65, 8).Value = 4055.6
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 4055.6
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 20278
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).Value = -26518

$ws.Cells.Item(122, 8).Value = 2019.0834
$ws.Cells.Item(122, 9).Value = 1943.1
$ws.Cells.Item(122, 10).Value = 2399
$ws.Cells.Item(122, 11).Value = 5829.299999999999
$ws.Cells.Item(122, 12).Value = 7197
$ws.Cells.Item(122, 13).Value = -3379.299999999999
$ws.Cells.Item(122, 14).Value = -12097

$ws.Cells.Item(126, 8).Value = 3000
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -6530
$ws.Cells.Item(126, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 2049.0715
$ws.Cells.Item(132, 10).Value = 2879
$ws.Cells.Item(132, 12).Value = 8637
$ws.Cells.Item(132, 14).Value = -13697
